{"js": "// Remove the merchant-assumptions bullet that claimed merchants will always\n// submit valid input (superseded now that input validation exists):\n// \"Merchants will provide valid input data when submitting payment requests.\"\nconst body = context.document.body;\n\nconst results = body.search(\n  \"Merchants will provide valid input data when submitting payment requests.\",\n  { matchCase: true, matchWholeWord: false }\n);\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\n    \"Target sentence not found: 'Merchants will provide valid input data when submitting payment requests.'\"\n  );\n}\n\n// The search hit is the run text; grab its enclosing paragraph and delete the\n// whole paragraph (bullet item) from the document body.\nconst hitRange = results.items[0];\nconst paragraph = hitRange.paragraphs.getFirst();\nparagraph.delete();\n\nawait context.sync();\n", "ps1": "# Remove the merchant-assumptions bullet that claimed merchants will always\n# submit valid input (superseded now that input validation exists):\n# \"Merchants will provide valid input data when submitting payment requests.\"\n$d = $word.ActiveDocument\n\n$searchText = \"Merchants will provide valid input data when submitting payment requests.\"\n\n$finder = $d.Content\n$found = $finder.Find.Execute($searchText)\n\nif (-not $found) {\n    throw \"Target sentence not found: '$searchText'\"\n}\n\n# $finder is now collapsed to the matched text; grab the whole enclosing\n# paragraph (the bulleted list item) so we can remove it completely,\n# including its paragraph mark.\n$target = $finder.Paragraphs.First.Range\n\n# Build a delete range that spans from the start of this paragraph through\n# the start of the following paragraph so the paragraph mark goes with it\n# (plain Range.Delete leaves an empty paragraph behind). Fall back to\n# extending backwards over the previous paragraph mark if this is the last\n# paragraph in the story.\n$nextPara = $target.Next(4, 1)  # wdParagraph unit\nif ($nextPara -ne $null) {\n    $deleteRange = $d.Range($target.Start, $nextPara.Start)\n} else {\n    $deleteRange = $d.Range($target.Start - 1, $target.End)\n}\n\n$deleteRange.Delete()\n"}
